$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 32 (pushes the existing row 32 and everything
# below it down by one, e.g. old row 32 -> new row 33, ..., old row 135 ->
# new row 136). This mirrors the diff: dimension grows from A1:T135 to
# A1:T136 and every row from 33 on now holds what used to be one row above.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value = 2
$ws.Cells.Item(32, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44623
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100109
$ws.Cells.Item(32, 8).Value = "Uva"
$ws.Cells.Item(32, 9).Value = 100109001
$ws.Cells.Item(32, 10).Value = "Uva"
$ws.Cells.Item(32, 11).Value = "Flame Seedless"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 340
$ws.Cells.Item(32, 14).Value = 7000
$ws.Cells.Item(32, 15).Value = 8000
$ws.Cells.Item(32, 16).Value = 7500
$ws.Cells.Item(32, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(32, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 19).Value = 417
$ws.Cells.Item(32, 20).Value = 18
